$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: LUX -> Luxemburg (code 29)
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "LUX"
$ws.Cells.Item(31, 3).Value = "Luxemburg"

# Row 32: CD -> Democratic Republic of the Congo (code 30)
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = "CD"
$ws.Cells.Item(32, 3).Value = "Democratic Republic of the Congo"

# Apply the same style (s="1") as used in column A of other rows
$ws.Range("A30").Copy()
$ws.Range("A31:A32").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
